$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4408552348613739
$ws.Range("B1").Value = 3.874288320541382
$ws.Range("C1").Value = 5.826854228973389
$ws.Range("D1").Value = 1.630766630172729
$ws.Range("E1").Value = 0.9789308309555054
